$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100 (shifts old rows 100..140 down to 101..141)
$ws.Rows(100).Insert()

# Populate the newly inserted row 100 with the new weekly record
$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 44510
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112032
$ws.Range("G100").Value = "Zapallo italiano"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 120
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 11000
$ws.Range("M100").Value = 10500
$ws.Range("N100").Value = "$/caja 60 unidades"
$ws.Range("O100").Value = "Región del Maule"
$ws.Range("P100").Value = 175
$ws.Range("Q100").Value = 60
$ws.Range("R100").Value = "Hortaliza"
